$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D1").Value = "MI_SA_ET15_2070"
$ws.Range("E1").Value = "MI_SA_PCDec_ET15_2070"
$ws.Range("F1").Value = "MI_SA_PCDD_ET15_2070"
$ws.Range("G1").Value = "MI_SA_ETLow_2070"
$ws.Range("H1").Value = "MI_SA_PCDec_ETLow_2070"
$ws.Range("I1").Value = "MI_SA_PCDD_ETLow_2070"
